$d = $word.ActiveDocument

# 1) Fix double space before opening quote and "sidplay" -> "display" typo
$d.Content.Find.Execute(
    "visual content of  " + [char]8220 + "Aganatiq" + [char]8221 + " game and it has some sections: General, player elements, heads up sidplay (HuD), antagonistic elements and global elements.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "visual content of " + [char]8220 + "Aganatiq" + [char]8221 + " game and it has some sections: General, player elements, heads up display (HuD), antagonistic elements and global elements.",
    2) | Out-Null

# 2) Remove the stray "R" (and the trailing period before it) after "...are used in this game."
$d.Content.Find.Execute(
    "are used in this game.R",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "are used in this game",
    2) | Out-Null

# 3) "Furthermore" -> "Moreover"
$d.Content.Find.Execute(
    "and paying money. Furthermore, in damage state,",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "and paying money. Moreover, in damage state,",
    2) | Out-Null

# 4) "can be shown" -> "can be indicated"; "damage to player" -> "damage player";
#    restructure the parenthetical so it closes after Aganatiq with a ")" before " in every hit."
$d.Content.Find.Execute(
    "can be shown. This ball can damage to player (taxi which is driven by " + [char]8216 + "Aganatiq" + [char]8217 + " in every hit.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "can be indicated. This ball can damage player (taxi which is driven by " + [char]8216 + "Aganatiq" + [char]8217 + ") in every hit.",
    2) | Out-Null

# 5) Fix typo "antaganostic" -> "antagonistic"
$d.Content.Find.Execute(
    "antaganostic elements",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "antagonistic elements",
    2) | Out-Null

# 6) "main character's car." -> "main character "Aganatiq"'s car."
$d.Content.Find.Execute(
    "overlap with main character" + [char]8217 + "s car.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "overlap with main character " + [char]8220 + "Aganatiq" + [char]8221 + [char]8217 + "s car.",
    2) | Out-Null

# 7) "arial" -> "Arial" (capitalisation)
$d.Content.Find.Execute(
    "Font Type of this game is arial font.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Font Type of this game is Arial font.",
    2) | Out-Null

# 8) Move the "_GoBack" bookmark from its old spot to right after
#    "...of taxi can be demonstrated." (collapsed bookmark at that point).
$rng = $d.Content
$rng.Find.Execute("of taxi can be demonstrated.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rng.End, $rng.End)
$insertPoint.InsertAfter([char]1)
$markRng = $d.Range($rng.End, $rng.End + 1)
$d.Bookmarks.Add("_GoBack", $markRng) | Out-Null
$delRng = $d.Range($rng.End, $rng.End + 1)
$delRng.Delete() | Out-Null

Write-Output "done"
